$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 3115
$ws.Range("I51").Value = 3600
$ws.Range("J51").Value = 2872.5
$ws.Range("K51").Value = 3600
$ws.Range("L51").Value = 2872.5
$ws.Range("M51").Value = -3116
$ws.Range("N51").Value = -3840.5
$ws.Range("H58").Value = 6016.8125
$ws.Range("I58").Value = 272.41666
$ws.Range("J58").Value = 23250
$ws.Range("K58").Value = 817.2499799999999
$ws.Range("L58").Value = 69750
$ws.Range("M58").Value = -667.2499799999999
$ws.Range("N58").Value = -70050
$ws.Range("H61").Value = 143.35
$ws.Range("I61").Value = 143.35
$ws.Range("K61").Value = 430.05
$ws.Range("M61").Value = -258.05
$ws.Range("H113").Value = 2662.3704
$ws.Range("I113").Value = 2054.5
$ws.Range("J113").Value = 3148.6667
$ws.Range("K113").Value = 2054.5
$ws.Range("L113").Value = 3148.6667
$ws.Range("M113").Value = 1199.5
$ws.Range("N113").Value = -9656.6667
$ws.Range("H132").Value = 2803077.5
$ws.Range("I132").Value = 3863118
$ws.Range("J132").Value = 1542
$ws.Range("K132").Value = 11589354
$ws.Range("L132").Value = 4626
$ws.Range("M132").Value = -11586824
$ws.Range("N132").Value = -9686
$ws.Range("H138").Value = 1718.0834
$ws.Range("I138").Value = 1383.8788
$ws.Range("J138").Value = 2126.5557
$ws.Range("K138").Value = 4151.636399999999
$ws.Range("L138").Value = 6379.6671
$ws.Range("M138").Value = 988.3636000000006
$ws.Range("N138").Value = -16659.6671
$ws.Range("H141").Value = 1474.7042
$ws.Range("I141").Value = 812.1091
$ws.Range("J141").Value = 3752.375
$ws.Range("K141").Value = 2436.3273
$ws.Range("L141").Value = 11257.125
$ws.Range("M141").Value = 2743.6727
$ws.Range("N141").Value = -21617.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1402.83
$ws.Range("I32").Value = 1295.6211
$ws.Range("J32").Value = 3439.8
$ws.Range("K32").Value = 1295.6211
$ws.Range("L32").Value = 3439.8
$ws.Range("M32").Value = -1008.6211
$ws.Range("N32").Value = -4013.8
$ws.Range("H61").Value = 1469.8667
$ws.Range("I61").Value = 928.1818
$ws.Range("K61").Value = 928.1818
$ws.Range("M61").Value = -716.1818
$ws.Range("H74").Value = 623.5769
$ws.Range("I74").Value = 522.9796
$ws.Range("K74").Value = 522.9796
$ws.Range("M74").Value = 351.0204
$ws.Range("H77").Value = 623.5769
$ws.Range("I77").Value = 522.9796
$ws.Range("K77").Value = 2614.898
$ws.Range("M77").Value = 1753.102
$ws.Range("H112").Value = 17346.75
$ws.Range("J112").Value = 17346.75
$ws.Range("L112").Value = 17346.75
$ws.Range("N112").Value = -20300.75
$ws.Range("H114").Value = 24037.6
$ws.Range("J114").Value = 24037.6
$ws.Range("L114").Value = 24037.6
$ws.Range("N114").Value = -32715.6
$ws.Range("H119").Value = 31714.285
$ws.Range("J119").Value = 31714.285
$ws.Range("L119").Value = 31714.285
$ws.Range("N119").Value = -41390.285
$ws.Range("H132").Value = 4342.4043
$ws.Range("I132").Value = 4812.1514
$ws.Range("J132").Value = 3235.1428
$ws.Range("K132").Value = 14436.4542
$ws.Range("L132").Value = 9705.428400000001
$ws.Range("M132").Value = -11906.4542
$ws.Range("N132").Value = -14765.4284
$ws.Range("H136").Value = 1469.8667
$ws.Range("I136").Value = 928.1818
$ws.Range("K136").Value = 2784.5454
$ws.Range("M136").Value = -234.5454

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H97").Value = 8880
$ws.Range("I97").Value = 6740.875
$ws.Range("J97").Value = 17436.5
$ws.Range("K97").Value = 6740.875
$ws.Range("L97").Value = 17436.5
$ws.Range("M97").Value = -5749.875
$ws.Range("N97").Value = -19418.5
$ws.Range("H134").Value = 21770.56
$ws.Range("I134").Value = 29253.945
$ws.Range("K134").Value = 87761.83499999999
$ws.Range("M134").Value = -85226.83499999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4067522
$ws.Range("I31").Value = 2070.96
$ws.Range("J31").Value = 10419789
$ws.Range("K31").Value = 2070.96
$ws.Range("L31").Value = 10419789
$ws.Range("M31").Value = -1775.96
$ws.Range("N31").Value = -10420379
$ws.Range("H34").Value = 4067522
$ws.Range("I34").Value = 2070.96
$ws.Range("J34").Value = 10419789
$ws.Range("K34").Value = 2070.96
$ws.Range("L34").Value = 10419789
$ws.Range("M34").Value = -1868.96
$ws.Range("N34").Value = -10420193
$ws.Range("H58").Value = 5376974
$ws.Range("I58").Value = 660.8542
$ws.Range("J58").Value = 23810046
$ws.Range("K58").Value = 660.8542
$ws.Range("L58").Value = 23810046
$ws.Range("M58").Value = -457.8542
$ws.Range("N58").Value = -23810452
$ws.Range("H132").Value = 1787.4117
$ws.Range("I132").Value = 1942.7059
$ws.Range("K132").Value = 5828.1177
$ws.Range("M132").Value = -3298.1177
$ws.Range("H136").Value = 5376974
$ws.Range("I136").Value = 660.8542
$ws.Range("J136").Value = 23810046
$ws.Range("K136").Value = 1982.5626
$ws.Range("L136").Value = 71430138
$ws.Range("M136").Value = 567.4374
$ws.Range("N136").Value = -71435238

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 987.6896400000001
$ws.Range("I5").Value = 322
$ws.Range("J5").Value = 1807
$ws.Range("K5").Value = 966
$ws.Range("L5").Value = 5421
$ws.Range("M5").Value = -854
$ws.Range("N5").Value = -5645
$ws.Range("H107").Value = 351
$ws.Range("I107").Value = 590
$ws.Range("J107").Value = 307.54544
$ws.Range("K107").Value = 1770
$ws.Range("L107").Value = 922.63632
$ws.Range("M107").Value = 150
$ws.Range("N107").Value = -4762.63632
$ws.Range("H131").Value = 1441992.9
$ws.Range("I131").Value = 5437.6816
$ws.Range("J131").Value = 1863382.5
$ws.Range("K131").Value = 16313.0448
$ws.Range("L131").Value = 5590147.5
$ws.Range("M131").Value = -11273.0448
$ws.Range("N131").Value = -5600227.5
$ws.Range("H132").Value = 2211.5
$ws.Range("I132").Value = 1279
$ws.Range("J132").Value = 2646.6667
$ws.Range("K132").Value = 11511
$ws.Range("L132").Value = 23820.0003
$ws.Range("M132").Value = -8981
$ws.Range("N132").Value = -28880.0003
$ws.Range("H135").Value = 987.6896400000001
$ws.Range("I135").Value = 322
$ws.Range("J135").Value = 1807
$ws.Range("K135").Value = 2898
$ws.Range("L135").Value = 16263
$ws.Range("M135").Value = -363
$ws.Range("N135").Value = -21333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 26481.707
$ws.Range("I132").Value = 35309.367
$ws.Range("J132").Value = 2406.2727
$ws.Range("K132").Value = 105928.101
$ws.Range("L132").Value = 7218.8181
$ws.Range("M132").Value = -103398.101
$ws.Range("N132").Value = -12278.8181

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4502
$ws.Range("I132").Value = 5794.846
$ws.Range("J132").Value = 1536.0588
$ws.Range("K132").Value = 17384.538
$ws.Range("L132").Value = 4608.1764
$ws.Range("M132").Value = -14854.538
$ws.Range("N132").Value = -9668.1764
$ws.Range("H136").Value = 1947.6951
$ws.Range("I136").Value = 1940.9855
$ws.Range("J136").Value = 1983.3077
$ws.Range("K136").Value = 5822.9565
$ws.Range("L136").Value = 5949.9231
$ws.Range("M136").Value = -3272.9565
$ws.Range("N136").Value = -11049.9231

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 27848.5
$ws.Range("J119").Value = 27848.5
$ws.Range("L119").Value = 27848.5
$ws.Range("N119").Value = -37524.5
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H132").Value = 1136.7021
$ws.Range("I132").Value = 1123.025
$ws.Range("J132").Value = 1214.8572
$ws.Range("K132").Value = 3369.075
$ws.Range("L132").Value = 3644.5716
$ws.Range("M132").Value = -839.0750000000003
$ws.Range("N132").Value = -8704.571599999999
$ws.Range("H136").Value = 2078.7683
$ws.Range("I136").Value = 2302.6323
$ws.Range("J136").Value = 991.4286
$ws.Range("K136").Value = 6907.896900000001
$ws.Range("L136").Value = 2974.2858
$ws.Range("M136").Value = -4357.896900000001
$ws.Range("N136").Value = -8074.2858
